$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{
    B = -0.999391737847945
    C = -168212921.728772
    D = 0.962855769528591
    E = 0.999391280599276
    F = -0.962861929841827
    G = 137191.16226133
    H = 7064317.99834326
    I = -3266764.51060286
    J = 0.963938592567792
    K = 0.999483086929724
    L = -0.962861929841827
    M = 56.8097427097398
    N = 2922.25889993025
    O = -3266764.51060286
  }
  3 = @{
    B = 0.960692208521904
    C = 191.766895909241
    D = -0.989801480892533
    E = -0.962466664394755
    F = 0.985316499870559
    G = -14.1145160291415
    H = -680.884385471807
    I = 3.96456314240996
    J = -0.941145086554243
    K = -0.83308458247653
    L = 0.985316499870559
    M = -0.055190831530792
    N = -2.42364630721163
    O = 3.96456314240996
  }
  4 = @{
    B = -0.999391742941703
    C = -168213116.810482
    D = 0.962861972874171
    E = 0.999391739254949
    F = -0.962862005314365
    G = 137204.700074537
    H = 7064972.81960954
    I = -3266768.53858099
    J = 0.962855643078381
    K = 0.999391020084426
    L = -0.962862005314365
    M = 55.898545045192
    N = 2878.35630414182
    O = -3266768.53858099
  }
  5 = @{
    B = 0.961809438901286
    C = 194.081710310473
    D = -0.990159378487734
    E = -0.964998854783382
    F = 0.985485299274315
    G = -14.5378132072954
    H = -702.895176981353
    I = 4.00844509920539
    J = -0.994071220285752
    K = -0.941231183566587
    L = 0.985485299274315
    M = -0.0888023354520817
    N = -4.17131493677733
    O = 4.00844509920539
  }
}

foreach ($row in $data.Keys) {
  foreach ($col in $data[$row].Keys) {
    $ws.Range("$col$row").Value = $data[$row][$col]
  }
}
